# Apply the Thu Sep 14 22:15:05 UTC 2023 cryptos-list refresh (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value that must stay text even when it looks numeric
# (e.g. "1.00", "0.0170") - set a temporary text NumberFormat so Excel does not
# coerce it to a number, then restore the default "Normal" style afterwards so
# the cell keeps its original (unstyled) appearance.
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '26.605.26'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '1.631.21'
$ws.Range('E4').Value = '  +0.02%  '
Set-TextCell 'D5' '212.91'
$ws.Range('E5').Value = '  -0.04%  '
Set-TextCell 'D6' '0.494'
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.78%  '
Set-TextCell 'D9' '0.0624'
$ws.Range('E9').Value = '  +1.50%  '
Set-TextCell 'D10' '18.97'
$ws.Range('E10').Value = '  +2.67%  '
Set-TextCell 'D11' '0.0841'
$ws.Range('E11').Value = '  +3.52%  '
$ws.Range('D12').Value = '1.858.53'
$ws.Range('E12').Value = '  +1.38%  '
$ws.Range('D13').Value = '1.621.80'
$ws.Range('E13').Value = '  +0.73%  '
Set-TextCell 'D14' '4.09'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').Value = '26.578.73'
$ws.Range('E16').Value = '  +1.14%  '
Set-TextCell 'D17' '63.02'
$ws.Range('E17').Value = '  +1.25%  '
$ws.Range('E18').Value = '  +1.35%  '
Set-TextCell 'D19' '209.35'
$ws.Range('E19').Value = '  +3.59%  '
$ws.Range('E21').Value = '  +0.66%  '
Set-TextCell 'D22' '9.44'
$ws.Range('E22').Value = '  +1.03%  '
Set-TextCell 'D23' '6.20'
$ws.Range('E23').Value = '  +2.94%  '
$ws.Range('E24').Value = '  +2.40%  '
Set-TextCell 'D25' '146.96'
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -0.38%  '
Set-TextCell 'D28' '6.85'
$ws.Range('E28').Value = '  +4.18%  '
Set-TextCell 'D29' '15.37'
$ws.Range('E29').Value = '  +0.83%  '
Set-TextCell 'D30' '0.0522'
$ws.Range('E30').Value = '  +5.18%  '
$ws.Range('E31').Value = '  -0.16%  '
Set-TextCell 'D32' '3.25'
$ws.Range('E32').Value = '  +1.60%  '
Set-TextCell 'D33' '2.95'
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('E34').Value = '  +0.66%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').Value = '1.166.02'
Set-TextCell 'D37' '0.0170'
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('E38').Value = '  +2.23%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('E40').Value = '  +1.44%  '
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('E42').Value = '  +1.04%  '
Set-TextCell 'D43' '5.38'
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').Value = '1.768.64'
$ws.Range('E44').Value = '  +1.41%  '
Set-TextCell 'D45' '92.14'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('E46').Value = '  +0.54%  '
Set-TextCell 'D47' '54.67'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D48' '0.0510'
$ws.Range('E48').Value = '  +0.66%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D49' '7.56'
$ws.Range('E49').Value = '  +4.00%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D50' '0.409'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextCell 'D51' '1.00'
$ws.Range('E51').Value = '  -0.13%  '
